$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D column (Price) holds numeric-looking text; force text type via a temporary
# "@" (Text) number format so Excel does not reinterpret values like "1.00" or
# "29.50" as numbers (which would drop formatting/precision), then clear the
# format again so the cell keeps its original style (no explicit s= attribute).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "65.786.84"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value2 = "  +0.58%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "2.673.76"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value2 = "  +0.99%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "601.88"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value2 = "  -0.56%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "157.06"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value2 = "  +0.70%  "

$ws.Range("E7").Value2 = "  -0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.616"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value2 = "  +5.05%  "

$ws.Range("E9").Value2 = "  +1.32%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "5.94"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value2 = "  +1.47%  "

$ws.Range("E11").Value2 = "  -0.19%  "

$ws.Range("E12").Value2 = "  -0.26%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "29.50"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value2 = "  -1.36%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "0.0000198"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value2 = "  +1.43%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "3.153.89"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value2 = "  +0.89%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "65.578.75"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value2 = "  +0.60%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "2.676.64"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value2 = "  +1.09%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "12.64"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value2 = "  -0.59%  "

$ws.Range("E19").Value2 = "  -1.22%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "7.57"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value2 = "  +1.33%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "352.30"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value2 = "  -1.65%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "1.00"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value2 = "  -0.19%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "69.89"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value2 = "  +0.35%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "0.0000111"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value2 = "  +5.81%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "9.81"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value2 = "  +4.19%  "

$ws.Range("E26").Value2 = "  -4.15%  "

$ws.Range("E27").Value2 = "  +1.28%  "

$ws.Range("E28").Value2 = "  -1.31%  "

$ws.Range("E29").Value2 = "  +0.39%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "544.99"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value2 = "  +4.56%  "

$ws.Range("E31").Value2 = "  -0.12%  "

$ws.Range("E32").Value2 = "  -0.63%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "1.78"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value2 = "  +0.20%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "6.58"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value2 = "  +3.90%  "

$ws.Range("E35").Value2 = "  -0.76%  "

$ws.Range("E36").Value2 = "  -1.87%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "20.42"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value2 = "  -1.17%  "

$ws.Range("E38").Value2 = "  +0.01%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "158.79"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value2 = "  -2.27%  "

$ws.Range("E40").Value2 = "  -0.81%  "

$ws.Range("E41").Value2 = "  +0.02%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "42.76"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value2 = "  +1.99%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "164.88"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value2 = "  -0.34%  "

$ws.Range("E44").Value2 = "  -0.58%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "0.0614"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value2 = "  +0.97%  "

$ws.Range("E46").Value2 = "  -1.06%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "23.25"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value2 = "  +1.58%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "0.645"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value2 = "  -0.84%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "0.0259"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value2 = "  -1.43%  "

$ws.Range("E51").Value2 = "  +3.30%  "
